$wb = $excel.ActiveWorkbook
Write-Host ($wb.Styles | Get-Member | Out-String)
